$d = $word.ActiveDocument

# --- 1) Mark every inline image that doesn't yet have NoProofing set as
#        NoProofing = True, matching Word's usual "picture run" formatting
#        (mirrors the other pre-existing images in the document, which
#        already carry <w:noProof/>). ---
$updated = 0
for ($i = 1; $i -le $d.InlineShapes.Count; $i++) {
    $ishp = $d.InlineShapes.Item($i)
    if (-not $ishp.Range.NoProofing) {
        $ishp.Range.NoProofing = 1
        $updated = $updated + 1
    }
}
Write-Output "NoProofing newly set on $updated inline image(s)"

# --- 2) Add a trailing space run to the final (empty) paragraph, matching
#        the Times New Roman / 24pt direct formatting used throughout the
#        document. We borrow formatting (and text) from an existing plain
#        run that is known to contain exactly a single space with no
#        rsid attributes, located via a unique anchor text right after it. ---
$searchRange = $d.Content
$found = $searchRange.Find.Execute("Entendendo a Web:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not locate anchor text to borrow space-run formatting from"
}
$srcRange = $d.Range($searchRange.Start - 1, $searchRange.Start)

$p = $d.Paragraphs.Last
$before_end = $p.Range.End
$p.Range.InsertAfter(" ")
$newRange = $d.Range($before_end - 1, $before_end)
$newRange.FormattedText = $srcRange.FormattedText
Write-Output "Appended formatted space run to final paragraph"
